$d = $word.ActiveDocument

# The document's last paragraph currently holds the "星期六..." diary entry
# and also carries a bookmark (_GoBack). Per the diff we need to:
#   1. Insert two brand-new paragraphs immediately BEFORE that last paragraph:
#        - "星期六，晴，今天天气很好，我和朋友一起出去玩了。又是开心的一天呢，心情也很好"
#        - "2023年3月16日 星期四"
#      Both using rFonts hint="eastAsia" on both the paragraph mark run
#      properties (pPr/rPr) and the text run properties (r/rPr), matching
#      the style used throughout the rest of the document.
#   2. Change the text of the (still) last paragraph (which keeps its
#      original rFonts hint="default" on pPr/rPr and hint="eastAsia" on the
#      run) from the old diary entry to "今天的作业很难，我很难过".

$last = $d.Paragraphs.Last

# --- Step 1: update the text of the last paragraph first, while its text is
# still unique in the document (avoids the upcoming insert, in step 2,
# creating a duplicate match for this Find) ---
$last.Range.Find.Execute("星期六，晴，今天天气很好，我和朋友一起出去玩了。又是开心的一天呢，心情也很好", $true, $false, $false, $false, $false, $true, 1, $false, "今天的作业很难，我很难过", 2)

# --- Step 2: insert the two new paragraphs before the last paragraph ---
# Build a raw OOXML fragment with the two desired paragraphs plus a trailing
# empty paragraph marker; InsertXML splits the destination paragraph at the
# insertion point, so the extra trailing <w:p/> keeps our two paragraphs
# intact and isolated from the original last paragraph. We remove that
# trailing helper paragraph immediately afterward.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$para1 = '<w:p ' + $wNs + '><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>星期六，晴，今天天气很好，我和朋友一起出去玩了。又是开心的一天呢，心情也很好</w:t></w:r></w:p>'
$para2 = '<w:p ' + $wNs + '><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>2023年3月16日 星期四</w:t></w:r></w:p>'
$paraSpacer = '<w:p ' + $wNs + '></w:p>'

$insPoint = $d.Range($last.Range.Start, $last.Range.Start)
$insPoint.InsertXML($para1 + $para2 + $paraSpacer)

# Remove the helper spacer paragraph that InsertXML left behind right before
# the original last paragraph.
$spacerPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$spacerPara.Range.Delete()
